$wb = $excel.ActiveWorkbook

# --- 1. Move the "Seats" worksheet tab so it sits right after "Leg_Schedule" ---
# (drags Seats from its old spot at the end of the table sheets up to position 5)
$seats = $wb.Worksheets.Item("Seats")
$legSchedule = $wb.Worksheets.Item("Leg_Schedule")
$seats.Move($null, $legSchedule)

# re-fetch a live reference to the worksheet post-move
$seats = $wb.Worksheets.Item("Seats")

# --- 2. Populate the Seats table with 15 seat rows (flight 1, legs 1-3, all available) ---
$seatData = @(
    @(1, 1, 1),
    @(2, 1, 1),
    @(3, 1, 1),
    @(4, 1, 1),
    @(5, 1, 1),
    @(6, 1, 2),
    @(7, 1, 2),
    @(8, 1, 2),
    @(9, 1, 2),
    @(10, 1, 2),
    @(11, 1, 3),
    @(12, 1, 3),
    @(13, 1, 3),
    @(14, 1, 3),
    @(15, 1, 3)
)

$row = 3
foreach ($r in $seatData) {
    $seats.Cells.Item($row, 1).Value = $r[0]   # A: seat_number
    $seats.Cells.Item($row, 4).Value = $r[1]   # D: flight_number
    $seats.Cells.Item($row, 5).Value = $r[2]   # E: leg_number
    $seats.Cells.Item($row, 6).Value = "Y"     # F: available
    $row++
}

# --- 3. Restore / update per-sheet selections that Excel remembered ---
$fares = $wb.Worksheets.Item("Fares")
$fares.Activate()
$fares.Range("D5").Select() | Out-Null

$legSchedule = $wb.Worksheets.Item("Leg_Schedule")
$legSchedule.Activate()
$legSchedule.Range("A3:B3").Select() | Out-Null

# --- 4. Seats becomes the active / selected tab, with C8 selected ---
$seats = $wb.Worksheets.Item("Seats")
$seats.Activate()
$seats.Range("C8").Select() | Out-Null

# --- 5. Scroll the workbook tab strip so "Flights" is the first visible tab ---
$excel.ActiveWindow.ScrollWorkbookTabs(2) | Out-Null
